$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Statement start balance date
$ws.Range("D5").Value = "KONTOSTAND AM 29.04.2024"

# Row 6
$ws.Range("B6").Value = "02.05."
$ws.Range("C6").Value = "03.05."
$ws.Range("D6").Value = "PAYPAL SMTKFA"
$ws.Range("E6").Value = "60,87-"

# Row 7
$ws.Range("B7").Value = "06.05."
$ws.Range("C7").Value = "07.05."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU GSPYJD"
$ws.Range("E7").Value = "147,78-"

# Row 8
$ws.Range("B8").Value = "07.05."
$ws.Range("C8").Value = "08.05."
$ws.Range("D8").Value = "PAYPAL YMJWZH"
$ws.Range("E8").Value = "18,27-"

# Row 9
$ws.Range("B9").Value = "09.05."
$ws.Range("C9").Value = "10.05."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 59613157"
$ws.Range("E9").Value = "86,39-"

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 12.05.2024"
$ws.Range("E12").Value = "313,31-"

# Next settlement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.05.2024"
